$d = $word.ActiveDocument

# Locate the three paragraphs that need to be replaced:
#  - "LP highlight yellow case_id when data is entered, DN highlight green when verified"
#  - the blank paragraph right after it
#  - the paragraph holding the _GoBack bookmark
$p28 = $d.Paragraphs.Item(28)
$p30 = $d.Paragraphs.Item(30)

$target = $d.Range($p28.Range.Start, $p30.Range.End)

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$xml = @"
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:rStyle w:val="e24kjd"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="e24kjd"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Case_id</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="e24kjd"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> in green = has been entered and verified</w:t>
  </w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:rStyle w:val="e24kjd"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="e24kjd"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Case_id</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rStyle w:val="e24kjd"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> in yellow = has been entered with questions/comments in the Comments column</w:t>
  </w:r>
</w:p>
<w:p $w>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:rPr>
      <w:rStyle w:val="e24kjd"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
<w:p $w/>
<w:p $w/>
"@

$target.InsertXML($xml)
